# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) figures to the
# "展览" and "全部类型" sheets, which hold identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G4").Value = 40

    $ws.Range("F5").Value = 112
    $ws.Range("G5").Value = 65

    $ws.Range("G6").Value = "不可售"

    $ws.Range("F9").Value = 603

    $ws.Range("F11").Value = 326

    $ws.Range("F13").Value = 389

    $ws.Range("F17").Value = 56

    $ws.Range("F19").Value = 104

    $ws.Range("F20").Value = 1023

    $ws.Range("F21").Value = 1426

    $ws.Range("F22").Value = 311

    $ws.Range("F27").Value = 46

    $ws.Range("F38").Value = 3812

    $ws.Range("F40").Value = 443

    $ws.Range("F41").Value = 219

    $ws.Range("F42").Value = 955

    $ws.Range("F43").Value = 60

    $ws.Range("F46").Value = 82
}
